$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76 (pushes existing rows 76-104 down to 77-105,
# which matches the content alignment seen in the target diff), then populate
# it with the new Berenjena record for Macroferia Regional de Talca.
$ws.Rows.Item(76).Insert()

$ws.Cells.Item(76, 1).Value = 5
$ws.Cells.Item(76, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(76, 3).Value = "Maule"
$ws.Cells.Item(76, 4).Value = 44609
$ws.Cells.Item(76, 5).Value = 7
$ws.Cells.Item(76, 6).Value = 100112001
$ws.Cells.Item(76, 7).Value = "Berenjena"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 150
$ws.Cells.Item(76, 11).Value = 7000
$ws.Cells.Item(76, 12).Value = 7000
$ws.Cells.Item(76, 13).Value = 7000
$ws.Cells.Item(76, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(76, 15).Value = "Región del Maule"
$ws.Cells.Item(76, 16).Value = 140
$ws.Cells.Item(76, 17).Value = 50
$ws.Cells.Item(76, 18).Value = "Hortaliza"
